$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B1").Value = 15
$ws.Range("C1").Value = 16
$ws.Range("D1").Value = 15
$ws.Range("E1").Value = 16

$ws.Range("B2").Value = 10.923144168135082
$ws.Range("C2").Value = 22.475702068216322
$ws.Range("D2").Value = 28.944811245278572
$ws.Range("E2").Value = 22.865305192092251

$ws.Range("B3").Value = 8.8236132846301416
$ws.Range("C3").Value = 13.472937611358049
$ws.Range("D3").Value = 39.914604100960105
$ws.Range("E3").Value = 15.5553783807012

$ws.Range("B1:E3").Select()
